$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Date line: "03. 09. 2019" -> "03. 12. 2019", with the "_GoBack" bookmark
#    relocated to sit right after "03. 12" (where the edit cursor landed).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("03. 09. 2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "03. 12. 2019", 2) | Out-Null

# The original "_GoBack" bookmark lives in the "PA 3 reflection" line; drop it
# so it can be re-created at the new location.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Re-create "_GoBack" right after "03. 12" in the date paragraph.
$datePara = $d.Paragraphs.Item(2)
$bmPos = $datePara.Range.Start + 6
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) "PA 3 " + "reflection" runs (previously split by the bookmark) collapse
#    into a single run once the bookmark is gone.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("PA 3 reflection", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "PA 3 reflection", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the trailing sentence from the first reflection paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This assignment seemed to be simpler than I thought as the most important code was given to us already. ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Replace the whole "One thing that threw me off..." paragraph with the
#    new tier-3 / Dijkstra reflection text.
# ---------------------------------------------------------------------------
$oldPara4 = "One thing that threw me off was that I forgot to fetch for the last change done on the PA3, " + `
    "so for one of my files, I was missing a connection and was getting a different result. " + `
    "It was producing the correct time for that map, but it wasn’t the right map. Other than that, " + `
    "that was my main issue. I tend to overlook the details I think wouldn’t be the issue, and in the end, " + `
    "it turns out to be issue. "

$newPara4 = "For tier 3, I got 38 minutes and 27 minutes when I ran deliveries 3.1/3.2 at first because I was computing " + `
    "Dijkstra’s algorithm(compute_shortest_path) at each delivery and then I searched for the next delivery in the " + `
    "map that was returned for computing the shortest path. Then once having new weights and new connections, I did " + `
    "the MST. I assumed it was that simple because my answers were right for the previous maps and deliveries."

$d.Content.Find.Execute($oldPara4, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newPara4, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Tweak wording in the "What really helped..." paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("made it easier to program because", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "made it possible to understand the program because", 2) | Out-Null
